$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.875.45'
$ws.Range("E2").Value = '  +6.73%  '
$ws.Range("D3").Value = '3.113.87'
$ws.Range("E3").Value = '  +3.90%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.01'
$ws.Range("E5").Value = '  +4.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.98'
$ws.Range("E6").Value = '  +5.16%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '3.104.44'
$ws.Range("E8").Value = '  +3.85%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.533'
$ws.Range("E9").Value = '  +2.26%  '
$ws.Range("E10").Value = '  +9.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.75'
$ws.Range("E11").Value = '  +10.74%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.470'
$ws.Range("E12").Value = '  +3.00%  '
$ws.Range("E13").Value = '  +6.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.63'
$ws.Range("E14").Value = '  +6.14%  '
$ws.Range("E15").Value = '  +0.88%  '
$ws.Range("D16").Value = '3.633.47'
$ws.Range("E16").Value = '  +4.05%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.30'
$ws.Range("E17").Value = '  -0.11%  '
$ws.Range("D18").Value = '3.115.80'
$ws.Range("E18").Value = '  +3.99%  '
$ws.Range("D19").Value = '62.790.15'
$ws.Range("E19").Value = '  +6.51%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '454.40'
$ws.Range("E20").Value = '  +6.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.11'
$ws.Range("E21").Value = '  +2.73%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.736'
$ws.Range("E22").Value = '  +1.83%  '
$ws.Range("E23").Value = '  +6.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.75'
$ws.Range("E24").Value = '  +3.40%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.42'
$ws.Range("E25").Value = '  +2.37%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.27'
$ws.Range("E27").Value = '  +3.93%  '
$ws.Range("E28").Value = '  +6.32%  '
$ws.Range("B29").Value = 'FirstDigitalUSD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.27'
$ws.Range("E30").Value = '  +6.39%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.85'
$ws.Range("E31").Value = '  +13.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.112'
$ws.Range("E32").Value = '  +13.47%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.16'
$ws.Range("E33").Value = '  +5.44%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.04'
$ws.Range("E34").Value = '  +4.36%  '
$ws.Range("D35").Value = '0.0₃0808'
$ws.Range("E35").Value = '  +6.98%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.10'
$ws.Range("E36").Value = '  +2.45%  '
$ws.Range("E37").Value = '  +7.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.06'
$ws.Range("E38").Value = '  +12.66%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '50.71'
$ws.Range("E39").Value = '  +4.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.81'
$ws.Range("E40").Value = '  +1.49%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '429.27'
$ws.Range("E41").Value = '  +7.95%  '
$ws.Range("D42").Value = '2.934.41'
$ws.Range("E42").Value = '  +6.63%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0373'
$ws.Range("E43").Value = '  +6.27%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.283'
$ws.Range("E44").Value = '  +12.63%  '
$ws.Range("E45").Value = '  +3.46%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.16'
$ws.Range("E46").Value = '  +8.43%  '
$ws.Range("B47").Value = 'Monero'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '125.18'
$ws.Range("E47").Value = '  +1.50%  '
$ws.Range("B48").Value = 'Arweave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '35.18'
$ws.Range("E48").Value = '  -0.52%  '
$ws.Range("B49").Value = 'USDe'
$ws.Range("C49").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.999'
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("E50").Value = '  +1.52%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.66'
$ws.Range("E51").Value = '  +5.67%  '
